$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header changes
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Column E rows 2-12 become the text "RLS" instead of numeric 1
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = "RLS"
}

# Tiny floating point recalculation differences in F, G, H columns
$ws.Range("F2").Value = 2.084208330992817
$ws.Range("G2").Value = 3.262163426006347
$ws.Range("H2").Value = 1.772608240358271

$ws.Range("F3").Value = 1.864381456141465
$ws.Range("G3").Value = 2.918094562769553
$ws.Range("H3").Value = 1.589002119133734

$ws.Range("F4").Value = 1.664860893028553
$ws.Range("G4").Value = 2.605808754271182
$ws.Range("H4").Value = 1.426433817933692

$ws.Range("F5").Value = 1.472634507535951
$ws.Range("G5").Value = 2.304939654506738
$ws.Range("H5").Value = 1.269610317792454

$ws.Range("F6").Value = 1.262047792143788
$ws.Range("G6").Value = 1.975333313941021
$ws.Range("H6").Value = 1.095677053477569

$ws.Range("F7").Value = 1.019300767019919
$ws.Range("G7").Value = 1.595390265371728
$ws.Range("H7").Value = 0.8811430688104318

$ws.Range("F8").Value = 0.8230228538669245
$ws.Range("G8").Value = 1.288179791207878
$ws.Range("H8").Value = 0.7042572388823858

$ws.Range("F9").Value = 0.7154428825532329
$ws.Range("G9").Value = 1.119797656575894
$ws.Range("H9").Value = 0.6021096348044288

$ws.Range("F10").Value = 0.6386071911451857
$ws.Range("G10").Value = 0.9995358868689067
$ws.Range("H10").Value = 0.5239262175316918

$ws.Range("F11").Value = 0.5919603517434724
$ws.Range("G11").Value = 0.9265251368530603
$ws.Range("H11").Value = 0.45575723369985

$ws.Range("F12").Value = 0.5837428749683309
$ws.Range("G12").Value = 0.9136632977598667
$ws.Range("H12").Value = 0.4432770401993973

$wb.Save()
